# "version final sin errores"
# - Bump the CodeSystem Version metadata value from 0.4.0 to 0.7.0
# - Remove the Jurisdiction ("Chile") row from the Metadata sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$lastRow = $ws.UsedRange.Rows.Count

$versionRow = -1
$jurisdictionRow = -1

for ($r = 1; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Text
    if ($label -eq "Version") {
        $versionRow = $r
    }
    if ($label -eq "Jurisdiction") {
        $jurisdictionRow = $r
    }
}

if ($versionRow -ne -1) {
    $ws.Cells.Item($versionRow, 2).Value = "0.7.0"
}

if ($jurisdictionRow -ne -1) {
    $ws.Rows.Item($jurisdictionRow).Delete()
}
